$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2,4) "54.231.82"
Set-TextValue $ws.Cells.Item(2,5) "  +1.15%  "
Set-TextValue $ws.Cells.Item(3,4) "2.268.54"
Set-TextValue $ws.Cells.Item(3,5) "  +2.58%  "
Set-TextValue $ws.Cells.Item(4,5) "  +0.04%  "
Set-TextValue $ws.Cells.Item(5,4) "495.44"
Set-TextValue $ws.Cells.Item(5,5) "  +1.76%  "
Set-TextValue $ws.Cells.Item(6,4) "128.01"
Set-TextValue $ws.Cells.Item(6,5) "  +2.11%  "
Set-TextValue $ws.Cells.Item(7,4) "0.998"
Set-TextValue $ws.Cells.Item(7,5) "  +0.31%  "
Set-TextValue $ws.Cells.Item(8,5) "  +1.03%  "
Set-TextValue $ws.Cells.Item(9,4) "0.0963"
Set-TextValue $ws.Cells.Item(9,5) "  +4.34%  "
Set-TextValue $ws.Cells.Item(10,5) "  +2.10%  "
Set-TextValue $ws.Cells.Item(11,5) "  +3.78%  "
Set-TextValue $ws.Cells.Item(12,4) "4.70"
Set-TextValue $ws.Cells.Item(12,5) "  +1.21%  "
Set-TextValue $ws.Cells.Item(13,4) "2.673.54"
Set-TextValue $ws.Cells.Item(13,5) "  +2.58%  "
Set-TextValue $ws.Cells.Item(14,4) "22.19"
Set-TextValue $ws.Cells.Item(14,5) "  +4.39%  "
Set-TextValue $ws.Cells.Item(15,4) "54.174.40"
Set-TextValue $ws.Cells.Item(15,5) "  +1.19%  "
Set-TextValue $ws.Cells.Item(16,5) "  +1.01%  "
Set-TextValue $ws.Cells.Item(17,4) "2.266.14"
Set-TextValue $ws.Cells.Item(17,5) "  +2.40%  "
Set-TextValue $ws.Cells.Item(18,4) "10.10"
Set-TextValue $ws.Cells.Item(18,5) "  +4.44%  "
Set-TextValue $ws.Cells.Item(19,4) "4.10"
Set-TextValue $ws.Cells.Item(19,5) "  +3.14%  "
Set-TextValue $ws.Cells.Item(20,4) "303.37"
Set-TextValue $ws.Cells.Item(20,5) "  +2.52%  "
Set-TextValue $ws.Cells.Item(21,4) "6.44"
Set-TextValue $ws.Cells.Item(21,5) "  +4.50%  "
Set-TextValue $ws.Cells.Item(22,5) "  +0.19%  "
Set-TextValue $ws.Cells.Item(23,5) "  -2.79%  "
Set-TextValue $ws.Cells.Item(24,4) "0.997"
Set-TextValue $ws.Cells.Item(24,5) "  +0.15%  "
Set-TextValue $ws.Cells.Item(25,4) "2.379.80"
Set-TextValue $ws.Cells.Item(25,5) "  +2.91%  "
Set-TextValue $ws.Cells.Item(26,5) "  +1.86%  "
Set-TextValue $ws.Cells.Item(27,5) "  +2.27%  "
Set-TextValue $ws.Cells.Item(28,4) "170.93"
Set-TextValue $ws.Cells.Item(28,5) "  +4.85%  "
Set-TextValue $ws.Cells.Item(29,5) "  +1.81%  "
Set-TextValue $ws.Cells.Item(30,4) "0.0₃0680"
Set-TextValue $ws.Cells.Item(30,5) "  +1.30%  "
Set-TextValue $ws.Cells.Item(31,5) "  +1.30%  "
Set-TextValue $ws.Cells.Item(32,5) "  +2.33%  "
Set-TextValue $ws.Cells.Item(33,5) "  +0.13%  "
Set-TextValue $ws.Cells.Item(34,4) "17.69"
Set-TextValue $ws.Cells.Item(34,5) "  +2.07%  "
Set-TextValue $ws.Cells.Item(35,5) "  +0.38%  "
Set-TextValue $ws.Cells.Item(36,4) "0.897"
Set-TextValue $ws.Cells.Item(36,5) "  +6.41%  "
Set-TextValue $ws.Cells.Item(37,5) "  +1.44%  "
Set-TextValue $ws.Cells.Item(38,4) "3.70"
Set-TextValue $ws.Cells.Item(38,5) "  +3.43%  "
Set-TextValue $ws.Cells.Item(39,4) "35.81"
Set-TextValue $ws.Cells.Item(39,5) "  +1.72%  "
Set-TextValue $ws.Cells.Item(40,5) "  +1.15%  "
Set-TextValue $ws.Cells.Item(41,4) "1.40"
Set-TextValue $ws.Cells.Item(41,5) "  +2.17%  "
Set-TextValue $ws.Cells.Item(42,5) "  +2.72%  "
Set-TextValue $ws.Cells.Item(43,4) "125.91"
Set-TextValue $ws.Cells.Item(43,5) "  -0.80%  "
Set-TextValue $ws.Cells.Item(44,4) "4.76"
Set-TextValue $ws.Cells.Item(44,5) "  -1.63%  "
Set-TextValue $ws.Cells.Item(45,4) "0.0898"
Set-TextValue $ws.Cells.Item(45,5) "  +1.91%  "
Set-TextValue $ws.Cells.Item(46,4) "0.0488"
Set-TextValue $ws.Cells.Item(46,5) "  +3.40%  "
Set-TextValue $ws.Cells.Item(47,4) "0.545"
Set-TextValue $ws.Cells.Item(47,5) "  +1.74%  "
Set-TextValue $ws.Cells.Item(48,4) "237.86"
Set-TextValue $ws.Cells.Item(48,5) "  +1.63%  "
Set-TextValue $ws.Cells.Item(49,5) "  +1.05%  "
Set-TextValue $ws.Cells.Item(50,4) "0.0206"
Set-TextValue $ws.Cells.Item(50,5) "  +2.39%  "
Set-TextValue $ws.Cells.Item(51,5) "  +0.99%  "
